$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row data (row 8), mirroring the layout/formatting of row 7
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(8, 1).Value = 42654.743680555555
$ws.Cells.Item(8, 2).Value = $false

$ws.Cells.Item(8, 3).Value = 9848.4
$ws.Cells.Item(8, 4).Value = 9903.36
$ws.Cells.Item(8, 5).Value = 308
$ws.Cells.Item(8, 6).Value = 304.58999599999999

$ws.Cells.Item(8, 7).Value = $false

$ws.Cells.Item(8, 8).Value = -1.1100000000000001

$ws.Cells.Item(8, 9).Value = $false

# Re-fit the columns to the new data, matching the original sheet's
# auto-fit column behavior ("bestFit"). The emulated ColumnWidth setter
# snaps to 1/6-character increments, so these are the closest achievable
# widths to the recalculated best-fit values.
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 7.3333333333333333
$ws.Columns.Item(3).ColumnWidth = 8
$ws.Columns.Item(4).ColumnWidth = 10.3333333333333333
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 10
$ws.Columns.Item(7).ColumnWidth = 9.5
$ws.Columns.Item(8).ColumnWidth = 13.8333333333333333
$ws.Columns.Item(9).ColumnWidth = 11
